# Applies the 2024-01-15 cryptos-list refresh (GitHub Actions job) to Sheet1.
# For every data row (2-51) column E (Volume 1h) is rewritten; most rows also
# get a new column D (Price). Price/volume are plain text cells in the source
# workbook, so numeric-looking prices (e.g. "312.00", "1.00") must be forced
# to Text before the write - otherwise the COM layer coerces them to real
# numbers and silently drops the significant trailing zeros/dots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price text (or $null if Price is unchanged), whether that Price
# string parses as a plain number, and the new Volume(1h) text.
$rowData = @{
    2 = @{ D='42.670.79'; DIsNumeric=$false; E='  -0.61%  ' }
    3 = @{ D='2.516.46'; DIsNumeric=$false; E='  -1.36%  ' }
    4 = @{ D=$null; DIsNumeric=$false; E='  +0.07%  ' }
    5 = @{ D='312.00'; DIsNumeric=$true; E='  +2.57%  ' }
    6 = @{ D='95.12'; DIsNumeric=$true; E='  -3.07%  ' }
    7 = @{ D='0.578'; DIsNumeric=$true; E='  +0.32%  ' }
    8 = @{ D=$null; DIsNumeric=$false; E='  +0.04%  ' }
    9 = @{ D='0.537'; DIsNumeric=$true; E='  -1.70%  ' }
    10 = @{ D='36.19'; DIsNumeric=$true; E='  -1.83%  ' }
    11 = @{ D=$null; DIsNumeric=$false; E='  -1.86%  ' }
    12 = @{ D='7.69'; DIsNumeric=$true; E='  -0.63%  ' }
    13 = @{ D=$null; DIsNumeric=$false; E='  -2.34%  ' }
    14 = @{ D='2.903.39'; DIsNumeric=$false; E='  -1.31%  ' }
    15 = @{ D='15.60'; DIsNumeric=$true; E='  +4.20%  ' }
    16 = @{ D='2.498.04'; DIsNumeric=$false; E='  -0.82%  ' }
    17 = @{ D='0.856'; DIsNumeric=$true; E='  -2.44%  ' }
    18 = @{ D='42.731.48'; DIsNumeric=$false; E='  -0.65%  ' }
    19 = @{ D='13.14'; DIsNumeric=$true; E='  -3.70%  ' }
    20 = @{ D='0.0₃0965'; DIsNumeric=$false; E='  -2.71%  ' }
    21 = @{ D=$null; DIsNumeric=$false; E='  -1.13%  ' }
    22 = @{ D='71.20'; DIsNumeric=$true; E='  -1.00%  ' }
    23 = @{ D='251.59'; DIsNumeric=$true; E='  -0.91%  ' }
    24 = @{ D='2.95'; DIsNumeric=$true; E='  -0.50%  ' }
    25 = @{ D=$null; DIsNumeric=$false; E='  -2.18%  ' }
    26 = @{ D='26.85'; DIsNumeric=$true; E='  -3.83%  ' }
    27 = @{ D='1.00'; DIsNumeric=$true; E='  +0.14%  ' }
    28 = @{ D=$null; DIsNumeric=$false; E='  +11.84%  ' }
    29 = @{ D='39.11'; DIsNumeric=$true; E='  +3.48%  ' }
    30 = @{ D=$null; DIsNumeric=$false; E='  -1.20%  ' }
    31 = @{ D=$null; DIsNumeric=$false; E='  -3.07%  ' }
    32 = @{ D='156.97'; DIsNumeric=$true; E='  -0.78%  ' }
    33 = @{ D='19.81'; DIsNumeric=$true; E='  +3.72%  ' }
    34 = @{ D='3.32'; DIsNumeric=$true; E='  +0.54%  ' }
    35 = @{ D=$null; DIsNumeric=$false; E='  -4.24%  ' }
    36 = @{ D=$null; DIsNumeric=$false; E='  -2.54%  ' }
    37 = @{ D='2.60'; DIsNumeric=$true; E='  -5.48%  ' }
    38 = @{ D=$null; DIsNumeric=$false; E='  -2.59%  ' }
    39 = @{ D='24.30'; DIsNumeric=$true; E='  -6.56%  ' }
    40 = @{ D=$null; DIsNumeric=$false; E='  -0.29%  ' }
    41 = @{ D='2.11'; DIsNumeric=$true; E='  -0.69%  ' }
    42 = @{ D='3.83'; DIsNumeric=$true; E='  -1.59%  ' }
    43 = @{ D='3.35'; DIsNumeric=$true; E='  -2.25%  ' }
    44 = @{ D='2.067.20'; DIsNumeric=$false; E='  -1.09%  ' }
    45 = @{ D='1.00'; DIsNumeric=$true; E='  +0.04%  ' }
    46 = @{ D='0.0301'; DIsNumeric=$true; E='  -1.59%  ' }
    47 = @{ D='86.61'; DIsNumeric=$true; E='  +0.02%  ' }
    48 = @{ D=$null; DIsNumeric=$false; E='  -1.85%  ' }
    49 = @{ D='2.757.42'; DIsNumeric=$false; E='  -1.51%  ' }
    50 = @{ D='73.86'; DIsNumeric=$true; E='  -1.19%  ' }
    51 = @{ D=$null; DIsNumeric=$false; E='  -0.66%  ' }
}

foreach ($row in $rowData.Keys) {
    $info = $rowData[$row]
    $ws.Range("E" + $row).Value = $info.E

    if ($null -ne $info.D) {
        $dCell = $ws.Range("D" + $row)
        if ($info.DIsNumeric) {
            # Numeric-looking text ("312.00", "1.00", ...): pin the cell to
            # Text format so Excel stores the literal string, then drop the
            # format override again so the cell keeps its original (default) style.
            $dCell.NumberFormat = "@"
            $dCell.Value = $info.D
            $dCell.ClearFormats()
        } else {
            # Already non-numeric text (thousand-separator dots, subscript digits, ...)
            # - a plain assignment keeps it as text.
            $dCell.Value = $info.D
        }
    }
}
